# Adds two new calculation columns (K/L) to the Lake Erie P-loading sheet:
#  - Q11: TP in western basin as 60% of 2013 (D9) total-P load
#  - K14/L14: new header labels for the western-basin TP/SRP columns
#  - K15/L15: TP (60% of 2019 D15) and SRP (27% of that TP) for western basin

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "TP in western basin as 60% from Maccoux" figure, placed at Q11,
# referencing the 2013 basin total (D9).
$ws.Range("Q11").Formula = "=0.6*D9"

# New header labels for the western-basin TP / SRP columns (row 14).
$ws.Range("K14").Value = "TP in western basin as 60% from Maccoux"
$ws.Range("L14").Value = "SRP in western basin as 26% from Maccoux"

# New western-basin TP (60% of basin total) and SRP (27% of that TP)
# figures for 2019 (row 15).
$ws.Range("K15").Formula = "=0.6*D15"
$ws.Range("L15").Formula = "=0.27*K15"

# Move/collapse the active selection to K15, matching the saved view state.
$ws.Range("K15").Select()
